$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 50331.637
$ws.Range("J17").Value = 47077.977
$ws.Range("L17").Value = 141233.931
$ws.Range("N17").Value = -141569.931

$ws.Range("H33").Value = 722.0454999999999
$ws.Range("I33").Value = 541.61536
$ws.Range("J33").Value = 982.6667
$ws.Range("K33").Value = 541.61536
$ws.Range("L33").Value = 982.6667
$ws.Range("M33").Value = -312.61536
$ws.Range("N33").Value = -1440.6667

$ws.Range("H43").Value = 4883.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 4883.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 4883.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -5021.5

$ws.Range("H53").Value = 295.1154
$ws.Range("I53").Value = 164.33333
$ws.Range("J53").Value = 364.35294
$ws.Range("K53").Value = 164.33333
$ws.Range("L53").Value = 364.35294
$ws.Range("M53").Value = 472.66667
$ws.Range("N53").Value = -1638.35294

$ws.Range("H58").Value = 2149.468
$ws.Range("I58").Value = 348.125
$ws.Range("J58").Value = 2518.9744
$ws.Range("K58").Value = 1044.375
$ws.Range("L58").Value = 7556.9232
$ws.Range("M58").Value = -894.375
$ws.Range("N58").Value = -7856.9232

$ws.Range("H129").Value = 946.3
$ws.Range("I129").Value = 1059.125
$ws.Range("J129").Value = 936.48914
$ws.Range("K129").Value = 3177.375
$ws.Range("L129").Value = 2809.46742
$ws.Range("M129").Value = 1822.625
$ws.Range("N129").Value = -12809.46742

$ws.Range("H133").Value = 39500
$ws.Range("J133").Value = 39500
$ws.Range("L133").Value = 39500
$ws.Range("N133").Value = -49620

$ws.Range("H135").Value = 1773.2778
$ws.Range("I135").Value = 1358.8667
$ws.Range("K135").Value = 12229.8003
$ws.Range("M135").Value = -9694.800300000001

$ws.Range("H138").Value = 3453.9285
$ws.Range("I138").Value = 2323.7273
$ws.Range("J138").Value = 3854.9678
$ws.Range("K138").Value = 6971.1819
$ws.Range("L138").Value = 11564.9034
$ws.Range("M138").Value = -1831.1819
$ws.Range("N138").Value = -21844.9034

$ws.Range("H141").Value = 516322.4
$ws.Range("I141").Value = 1526.4117
$ws.Range("J141").Value = 1766541.2
$ws.Range("K141").Value = 4579.2351
$ws.Range("L141").Value = 5299623.6
$ws.Range("M141").Value = 600.7649000000001
$ws.Range("N141").Value = -5309983.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2626.4443
$ws.Range("I122").Value = 1909.2
$ws.Range("K122").Value = 5727.6
$ws.Range("M122").Value = -3277.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2588.4211
$ws.Range("I134").Value = 1698.75
$ws.Range("K134").Value = 5096.25
$ws.Range("M134").Value = -2561.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2856
$ws.Range("I31").Value = 2335.7568
$ws.Range("J31").Value = 3818.45
$ws.Range("K31").Value = 2335.7568
$ws.Range("L31").Value = 3818.45
$ws.Range("M31").Value = -2040.7568
$ws.Range("N31").Value = -4408.45

$ws.Range("H34").Value = 2856
$ws.Range("I34").Value = 2335.7568
$ws.Range("J34").Value = 3818.45
$ws.Range("K34").Value = 2335.7568
$ws.Range("L34").Value = 3818.45
$ws.Range("M34").Value = -2133.7568
$ws.Range("N34").Value = -4222.45

$ws.Range("H74").Value = 16796.416
$ws.Range("J74").Value = 16796.416
$ws.Range("L74").Value = 16796.416
$ws.Range("N74").Value = -18544.416

$ws.Range("H77").Value = 16796.416
$ws.Range("J77").Value = 16796.416
$ws.Range("L77").Value = 50389.24800000001
$ws.Range("N77").Value = -59125.24800000001

$ws.Range("H122").Value = 3423.8235
$ws.Range("I122").Value = 3053.8
$ws.Range("J122").Value = 3952.4285
$ws.Range("K122").Value = 9161.400000000001
$ws.Range("L122").Value = 11857.2855
$ws.Range("M122").Value = -6711.400000000001
$ws.Range("N122").Value = -16757.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4866.909
$ws.Range("I3").Value = 4976.6665
$ws.Range("J3").Value = 4825.75
$ws.Range("K3").Value = 14929.9995
$ws.Range("L3").Value = 14477.25
$ws.Range("M3").Value = -14817.9995
$ws.Range("N3").Value = -14701.25

$ws.Range("H87").Value = 9492
$ws.Range("J87").Value = 15816.667
$ws.Range("L87").Value = 47450.001
$ws.Range("N87").Value = -49946.001

$ws.Range("H90").Value = 9492
$ws.Range("J90").Value = 15816.667
$ws.Range("L90").Value = 142350.003
$ws.Range("N90").Value = -154830.003

$ws.Range("H130").Value = 2055.3333
$ws.Range("J130").Value = 2428.2856
$ws.Range("L130").Value = 7284.8568
$ws.Range("N130").Value = -17324.8568

$ws.Range("H131").Value = 1542.48
$ws.Range("J131").Value = 1137.4783
$ws.Range("L131").Value = 3412.4349
$ws.Range("N131").Value = -13492.4349

$ws.Range("H133").Value = 4454.1333
$ws.Range("I133").Value = 4855.25
$ws.Range("K133").Value = 14565.75
$ws.Range("M133").Value = -9505.75

$ws.Range("H134").Value = 3099.6924
$ws.Range("I134").Value = 2164.75
$ws.Range("J134").Value = 3515.2222
$ws.Range("K134").Value = 6494.25
$ws.Range("L134").Value = 10545.6666
$ws.Range("M134").Value = -1424.25
$ws.Range("N134").Value = -20685.6666

$ws.Range("H136").Value = 2377.4
$ws.Range("I136").Value = 1458.6
$ws.Range("J136").Value = 3755.6
$ws.Range("K136").Value = 4375.799999999999
$ws.Range("L136").Value = 11266.8
$ws.Range("M136").Value = 724.2000000000007
$ws.Range("N136").Value = -21466.8

$ws.Range("H137").Value = 2522.476
$ws.Range("I137").Value = 1784
$ws.Range("J137").Value = 3193.818
$ws.Range("K137").Value = 5352
$ws.Range("L137").Value = 9581.454000000002
$ws.Range("M137").Value = -252
$ws.Range("N137").Value = -19781.454

$ws.Range("H138").Value = 1483.9166
$ws.Range("I138").Value = 839.8889
$ws.Range("K138").Value = 2519.6667
$ws.Range("M138").Value = 2620.3333

$ws.Range("H139").Value = 6758874
$ws.Range("I139").Value = 10418360
$ws.Range("J139").Value = 2900.1538
$ws.Range("K139").Value = 31255080
$ws.Range("L139").Value = 8700.4614
$ws.Range("M139").Value = -31249940
$ws.Range("N139").Value = -18980.4614

$ws.Range("H140").Value = 6669743
$ws.Range("I140").Value = 16667411
$ws.Range("K140").Value = 50002233
$ws.Range("M140").Value = -49997053

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 70007
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 70007
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 70007
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -70353

$ws.Range("H30").Value = 70007
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 70007
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 70007
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -70217

$ws.Range("H47").Value = 50031
$ws.Range("J47").Value = 50031
$ws.Range("L47").Value = 50031
$ws.Range("N47").Value = -51167

$ws.Range("H52").Value = 3676.6667
$ws.Range("I52").Value = 1030
$ws.Range("K52").Value = 1030
$ws.Range("M52").Value = -771

$ws.Range("H126").Value = 3428.7646
$ws.Range("I126").Value = 2127
$ws.Range("J126").Value = 4340
$ws.Range("K126").Value = 6381
$ws.Range("L126").Value = 13020
$ws.Range("M126").Value = -3911
$ws.Range("N126").Value = -17960

$ws.Range("H132").Value = 3316.1462
$ws.Range("I132").Value = 3163.818
$ws.Range("K132").Value = 9491.454000000002
$ws.Range("M132").Value = -6961.454000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 146030.42
$ws.Range("I23").Value = 202041.2
$ws.Range("J23").Value = 6003.5
$ws.Range("K23").Value = 202041.2
$ws.Range("L23").Value = 6003.5
$ws.Range("M23").Value = -201811.2
$ws.Range("N23").Value = -6463.5

$ws.Range("H33").Value = 54512.75
$ws.Range("J33").Value = 54512.75
$ws.Range("L33").Value = 54512.75
$ws.Range("N33").Value = -55092.75

$ws.Range("H46").Value = 2568.1
$ws.Range("I46").Value = 595.25
$ws.Range("J46").Value = 3883.3333
$ws.Range("K46").Value = 595.25
$ws.Range("L46").Value = 3883.3333
$ws.Range("M46").Value = -407.25
$ws.Range("N46").Value = -4259.3333

$ws.Range("H122").Value = 3704.9443
$ws.Range("I122").Value = 2700
$ws.Range("J122").Value = 4961.125
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 14883.375
$ws.Range("M122").Value = -5650
$ws.Range("N122").Value = -19783.375

$ws.Range("H136").Value = 3030.2144
$ws.Range("I136").Value = 2633.9092
$ws.Range("K136").Value = 7901.7276
$ws.Range("M136").Value = -5351.7276

$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1020.4
$ws.Range("I81").Value = 1020.4
$ws.Range("K81").Value = 2040.8
$ws.Range("M81").Value = -979.8

$ws.Range("H84").Value = 1020.4
$ws.Range("I84").Value = 1020.4
$ws.Range("K84").Value = 10204
$ws.Range("M84").Value = -4900

$ws.Range("H122").Value = 836691.25
$ws.Range("I122").Value = 1113532.9
$ws.Range("J122").Value = 6166.3335
$ws.Range("K122").Value = 3340598.7
$ws.Range("L122").Value = 18499.0005
$ws.Range("M122").Value = -3338148.7
$ws.Range("N122").Value = -23399.0005

$ws.Range("H140").Value = 33248.383
$ws.Range("I140").Value = 15200
$ws.Range("J140").Value = 36529.91
$ws.Range("K140").Value = 15200
$ws.Range("L140").Value = 36529.91
$ws.Range("M140").Value = -10020
$ws.Range("N140").Value = -46889.91
